{"js": "// Locate the target paragraph: \"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 b \u043c\u043e\u0436\u043d\u043e \u0441\u043a\u0430\u0437\u0430\u0442\u044c, \u0447\u0442\u043e\n// \u043f\u0440\u0438 \u0443\u0432\u0435\u043b\u0438\u0447\u0435\u043d\u0438\u0438 \u0441\u0443\u043c\u043c\u044b \u0447\u0435\u043a\u0430 \u0432\" -> rewritten into a forecast sentence about tip\n// amounts. We search the body's paragraphs for the unique opening text and\n// edit that paragraph's runs directly so unrelated paragraphs (e.g. the very\n// similar-looking sentence about the correlation coefficient a few lines\n// above) are left untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst matches = body.search(\"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 b \u043c\u043e\u0436\u043d\u043e \u0441\u043a\u0430\u0437\u0430\u0442\u044c\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst targetRange = matches.items[0];\nconst targetParagraphs = targetRange.paragraphs;\ntargetParagraphs.load(\"items\");\nawait context.sync();\nconst p = targetParagraphs.items[0];\n\n// 1) Remove the standalone \"b\" run (it carries an explicit en-US language\n//    override that should not survive once the run becomes Russian text).\nlet res = p.search(\"b\", { matchCase: true });\nawait context.sync();\nres.items[0].delete();\nawait context.sync();\n\n// 2) Insert the new forecast text right after \"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \".\n//    Anchoring the insert on that neighboring run means the new text picks\n//    up its formatting (Times New Roman / no language override) instead of\n//    the one that used to belong to the deleted \"b\" run.\nres = p.search(\"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \", { matchCase: true });\nawait context.sync();\nres.items[0].insertText(\n  \"\u043f\u0440\u044f\u043c\u043e\u0439 \u0440\u0435\u0433\u0440\u0435\u0441\u0441\u0438\u0438 \u0438 \u0440\u0430\u0437\u0431\u0440\u043e\u0441\u0443 \u0442\u043e\u0447\u0435\u043a \u043c\u043e\u0436\u043d\u043e \u0441\u0434\u0435\u043b\u0430\u0442\u044c \u043f\u0440\u043e\u0433\u043d\u043e\u0437, \u0447\u0442\u043e \u043f\u0440\u0438 \u0441\u0442\u043e\u0438\u043c\u043e\u0441\u0442\u0438 \u0447\u0435\u043a\u0430 60\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 3) Trim \"\u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \" so the lead-in now simply reads \"\u041f\u043e \".\nres = p.search(\"\u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \", { matchCase: true });\nawait context.sync();\nres.items[0].delete();\nawait context.sync();\n\n// 4) Replace the old tail (\"60 \u043c\u043e\u0436\u043d\u043e \u0441\u043a\u0430\u0437\u0430\u0442\u044c, \u0447\u0442\u043e \u043f\u0440\u0438 \u0443\u0432\u0435\u043b\u0438\u0447\u0435\u043d\u0438\u0438 \u0441\u0443\u043c\u043c\u044b \u0447\u0435\u043a\u0430 \u0432\")\n//    with the new one, inserting the \"$\" currency markers and the tip range.\nres = p.search(\"60 \u043c\u043e\u0436\u043d\u043e \u0441\u043a\u0430\u0437\u0430\u0442\u044c, \u0447\u0442\u043e \u043f\u0440\u0438 \u0443\u0432\u0435\u043b\u0438\u0447\u0435\u043d\u0438\u0438 \u0441\u0443\u043c\u043c\u044b \u0447\u0435\u043a\u0430 \u0432\", { matchCase: true });\nawait context.sync();\nres.items[0].insertText(\n  \"60$ \u0440\u0430\u0437\u043c\u0435\u0440 \u0447\u0430\u0435\u0432\u044b\u0445 \u0431\u0443\u0434\u0435\u0442 \u043d\u0430\u0445\u043e\u0434\u0438\u0442\u044c\u0441\u044f \u0432 \u0434\u0438\u0430\u043f\u0430\u0437\u043e\u043d\u0435 2-10$\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Rewrites the final bullet of the \"\u0412\u044b\u0432\u043e\u0434\u044b\" section from\n#   \"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 b \u043c\u043e\u0436\u043d\u043e \u0441\u043a\u0430\u0437\u0430\u0442\u044c, \u0447\u0442\u043e \u043f\u0440\u0438 \u0443\u0432\u0435\u043b\u0438\u0447\u0435\u043d\u0438\u0438 \u0441\u0443\u043c\u043c\u044b \u0447\u0435\u043a\u0430 \u0432\"\n# into the forecast sentence about expected tip amounts. The very similar\n# sentence earlier in the document (\"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \u043a\u043e\u0440\u0440\u0435\u043b\u044f\u0446\u0438\u0438 \u0438\n# \u0448\u043a\u0430\u043b\u0435 \u0427\u0435\u0434\u0434\u043e\u043a\u0430 ...\") must stay untouched, so every Find is scoped to the\n# target paragraph only.\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph via a phrase that only exists there.\n$anchor = $d.Content\n$anchorFind = $anchor.Find\n$anchorFind.ClearFormatting()\n$anchorFind.Text = \"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 b \u043c\u043e\u0436\u043d\u043e \u0441\u043a\u0430\u0437\u0430\u0442\u044c\"\n$anchorFind.MatchCase = $true\n$anchorFind.Execute() | Out-Null\n\n$para = $anchor.Duplicate\n$para.Expand(4) | Out-Null   # wdParagraph\n\n# Word's COM Range objects don't auto-grow here after edits made through a\n# different (duplicated) Range, so re-derive the paragraph's extent from its\n# (stable) start after every mutation.\nfunction Refresh-Paragraph($range) {\n    $range.Collapse(1) | Out-Null  # wdCollapseStart\n    $range.Expand(4) | Out-Null    # wdParagraph\n}\n\n# 1) Delete the standalone \"b\" run. It carries an explicit en-US language\n#    override that should not survive once the run becomes Russian text;\n#    deleting it merges its neighbors (which share identical formatting\n#    without the language override).\n$bScope = $para.Duplicate\n$bFind = $bScope.Find\n$bFind.ClearFormatting()\n$bFind.Text = \"b\"\n$bFind.MatchCase = $true\n$bFind.Execute() | Out-Null\n$bScope.Delete()\nRefresh-Paragraph $para\n\n# 2) Insert the new forecast text right after \"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \".\n#    Anchoring on that run means the inserted text inherits its formatting\n#    (Times New Roman, no language override).\n$anchorScope = $para.Duplicate\n$anchorScopeFind = $anchorScope.Find\n$anchorScopeFind.ClearFormatting()\n$anchorScopeFind.Text = \"\u041f\u043e \u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \"\n$anchorScopeFind.MatchCase = $true\n$anchorScopeFind.Execute() | Out-Null\n$anchorScope.InsertAfter(\"\u043f\u0440\u044f\u043c\u043e\u0439 \u0440\u0435\u0433\u0440\u0435\u0441\u0441\u0438\u0438 \u0438 \u0440\u0430\u0437\u0431\u0440\u043e\u0441\u0443 \u0442\u043e\u0447\u0435\u043a \u043c\u043e\u0436\u043d\u043e \u0441\u0434\u0435\u043b\u0430\u0442\u044c \u043f\u0440\u043e\u0433\u043d\u043e\u0437, \u0447\u0442\u043e \u043f\u0440\u0438 \u0441\u0442\u043e\u0438\u043c\u043e\u0441\u0442\u0438 \u0447\u0435\u043a\u0430 60\")\nRefresh-Paragraph $para\n\n# 3) Trim \"\u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \" so the lead-in now simply reads \"\u041f\u043e \".\n$trimScope = $para.Duplicate\n$trimFind = $trimScope.Find\n$trimFind.ClearFormatting()\n$trimFind.Text = \"\u0432\u0435\u043b\u0438\u0447\u0438\u043d\u0435 \u043a\u043e\u044d\u0444\u0444\u0438\u0446\u0438\u0435\u043d\u0442\u0430 \"\n$trimFind.MatchCase = $true\n$trimFind.Execute() | Out-Null\n$trimScope.Delete()\nRefresh-Paragraph $para\n\n# 4) Replace the old tail with the new one, adding the \"$\" currency markers\n#    and the tip-amount range.\n$tailScope = $para.Duplicate\n$tailFind = $tailScope.Find\n$tailFind.ClearFormatting()\n$tailFind.Text = \"60 \u043c\u043e\u0436\u043d\u043e \u0441\u043a\u0430\u0437\u0430\u0442\u044c, \u0447\u0442\u043e \u043f\u0440\u0438 \u0443\u0432\u0435\u043b\u0438\u0447\u0435\u043d\u0438\u0438 \u0441\u0443\u043c\u043c\u044b \u0447\u0435\u043a\u0430 \u0432\"\n$tailFind.MatchCase = $true\n$tailFind.Execute() | Out-Null\n$tailScope.Text = \"60`$ \u0440\u0430\u0437\u043c\u0435\u0440 \u0447\u0430\u0435\u0432\u044b\u0445 \u0431\u0443\u0434\u0435\u0442 \u043d\u0430\u0445\u043e\u0434\u0438\u0442\u044c\u0441\u044f \u0432 \u0434\u0438\u0430\u043f\u0430\u0437\u043e\u043d\u0435 2-10`$\"\nRefresh-Paragraph $para\n"}
